# Weekly fruit/vegetable price update: a new record was reported for this
# market/product combination, so a new row is inserted right after the
# existing row 15 (becoming the new row 16), pushing every subsequent row
# (old 16..65) down by one (new 17..66).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 16, shifting rows 16-65 down to 17-66.
$ws.Rows(16).Insert()

# Populate the newly inserted row 16 with the new record's data.
$ws.Cells.Item(16, 1).Value = 3
$ws.Cells.Item(16, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(16, 3).Value = "Coquimbo"
$ws.Cells.Item(16, 4).Value = 45012
$ws.Cells.Item(16, 5).Value = 5
$ws.Cells.Item(16, 6).Value = "Fruta"
$ws.Cells.Item(16, 7).Value = 100108
$ws.Cells.Item(16, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(16, 9).Value = 100108004
$ws.Cells.Item(16, 10).Value = "Papaya"
$ws.Cells.Item(16, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(16, 12).Value = "Primera"
$ws.Cells.Item(16, 13).Value = 60
$ws.Cells.Item(16, 14).Value = 20000
$ws.Cells.Item(16, 15).Value = 20000
$ws.Cells.Item(16, 16).Value = 20000
$ws.Cells.Item(16, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(16, 18).Value = "Provincia del Elquí"
$ws.Cells.Item(16, 19).Value = 2000
$ws.Cells.Item(16, 20).Value = 10
